$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 12060.8
$ws.Range("J10").Value = 12060.8
$ws.Range("L10").Value = 12060.8
$ws.Range("N10").Value = -12646.8

$ws.Range("H15").Value = 999.44183
$ws.Range("I15").Value = 999.44183
$ws.Range("K15").Value = 2998.32549
$ws.Range("M15").Value = -2829.32549

$ws.Range("H28").Value = 42524.047
$ws.Range("I28").Value = 397.41666
$ws.Range("K28").Value = 397.41666
$ws.Range("M28").Value = 87.58334000000002

$ws.Range("H40").Value = 14426.667
$ws.Range("I40").Value = 1999
$ws.Range("J40").Value = 15314.357
$ws.Range("K40").Value = 1999
$ws.Range("L40").Value = 15314.357
$ws.Range("M40").Value = -1824
$ws.Range("N40").Value = -15664.357

$ws.Range("H137").Value = 339053.88
$ws.Range("I137").Value = 1647.125
$ws.Range("J137").Value = 765251.9
$ws.Range("K137").Value = 4941.375
$ws.Range("L137").Value = 2295755.7
$ws.Range("M137").Value = -2391.375
$ws.Range("N137").Value = -2300855.7

$ws.Range("H138").Value = 1388
$ws.Range("I138").Value = 859.625
$ws.Range("J138").Value = 2444.75
$ws.Range("K138").Value = 2578.875
$ws.Range("L138").Value = 7334.25
$ws.Range("M138").Value = 2561.125
$ws.Range("N138").Value = -17614.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10319.712
$ws.Range("I32").Value = 8874.918
$ws.Range("K32").Value = 8874.918
$ws.Range("M32").Value = -8587.918

$ws.Range("H61").Value = 1481.7727
$ws.Range("I61").Value = 1202.8422
$ws.Range("K61").Value = 1202.8422
$ws.Range("M61").Value = -990.8422

$ws.Range("H98").Value = 49937.8
$ws.Range("I98").Value = 32200
$ws.Range("K98").Value = 32200
$ws.Range("M98").Value = -29205

$ws.Range("H110").Value = 2213.3333
$ws.Range("I110").Value = 1695
$ws.Range("J110").Value = 3250
$ws.Range("K110").Value = 1695
$ws.Range("L110").Value = 3250
$ws.Range("M110").Value = 350
$ws.Range("N110").Value = -7340

$ws.Range("H136").Value = 1481.7727
$ws.Range("I136").Value = 1202.8422
$ws.Range("K136").Value = 3608.5266
$ws.Range("M136").Value = -1058.5266

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 12464.167
$ws.Range("J100").Value = 12464.167
$ws.Range("L100").Value = 12464.167
$ws.Range("N100").Value = -14628.167

$ws.Range("H140").Value = 43486.8
$ws.Range("J140").Value = 43486.8
$ws.Range("L140").Value = 43486.8
$ws.Range("N140").Value = -53846.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22328.4
$ws.Range("J28").Value = 22328.4
$ws.Range("L28").Value = 22328.4
$ws.Range("N28").Value = -22818.4

$ws.Range("H31").Value = 2474.6
$ws.Range("I31").Value = 1441
$ws.Range("J31").Value = 3428.6924
$ws.Range("K31").Value = 1441
$ws.Range("L31").Value = 3428.6924
$ws.Range("M31").Value = -1146
$ws.Range("N31").Value = -4018.6924

$ws.Range("H34").Value = 2474.6
$ws.Range("I34").Value = 1441
$ws.Range("J34").Value = 3428.6924
$ws.Range("K34").Value = 1441
$ws.Range("L34").Value = 3428.6924
$ws.Range("M34").Value = -1239
$ws.Range("N34").Value = -3832.6924

$ws.Range("H105").Value = 114700.4
$ws.Range("J105").Value = 4749.8335
$ws.Range("L105").Value = 4749.8335
$ws.Range("N105").Value = -8243.833500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 66775
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 77046.69500000001
$ws.Range("K23").Value = 27
$ws.Range("L23").Value = 231140.085
$ws.Range("M23").Value = 208
$ws.Range("N23").Value = -231610.085

$ws.Range("H92").Value = 409.8
$ws.Range("I92").Value = 409.8
$ws.Range("K92").Value = 1229.4
$ws.Range("M92").Value = 18.59999999999991

$ws.Range("H113").Value = 67612.664
$ws.Range("I113").Value = 1012.25
$ws.Range("J113").Value = 91831
$ws.Range("K113").Value = 3036.75
$ws.Range("L113").Value = 275493
$ws.Range("M113").Value = -866.75
$ws.Range("N113").Value = -279833

$ws.Range("H132").Value = 6867.2
$ws.Range("J132").Value = 8134.25
$ws.Range("L132").Value = 73208.25
$ws.Range("N132").Value = -78268.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 3824.9
$ws.Range("I22").Value = 6466.6665
$ws.Range("J22").Value = 2692.7144
$ws.Range("K22").Value = 6466.6665
$ws.Range("L22").Value = 2692.7144
$ws.Range("M22").Value = -5937.6665
$ws.Range("N22").Value = -3750.7144

$ws.Range("H93").Value = 19579.445
$ws.Range("J93").Value = 19579.445
$ws.Range("L93").Value = 19579.445
$ws.Range("N93").Value = -23323.445

$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164

$ws.Range("H113").Value = 160867.28
$ws.Range("I113").Value = 279205.25
$ws.Range("J113").Value = 3083.3333
$ws.Range("K113").Value = 279205.25
$ws.Range("L113").Value = 3083.3333
$ws.Range("M113").Value = -277035.25
$ws.Range("N113").Value = -7423.3333

$ws.Range("H132").Value = 3444.7856
$ws.Range("I132").Value = 2313.5789
$ws.Range("K132").Value = 6940.736699999999
$ws.Range("M132").Value = -4410.736699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2680.3103
$ws.Range("I55").Value = 1535.2354
$ws.Range("J55").Value = 4302.5
$ws.Range("K55").Value = 1535.2354
$ws.Range("L55").Value = 4302.5
$ws.Range("M55").Value = -1362.2354
$ws.Range("N55").Value = -4648.5

$ws.Range("H136").Value = 2272.879
$ws.Range("I136").Value = 2499.2273
$ws.Range("J136").Value = 1820.1818
$ws.Range("K136").Value = 7497.6819
$ws.Range("L136").Value = 5460.5454
$ws.Range("M136").Value = -4947.6819
$ws.Range("N136").Value = -10560.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 43933.332
$ws.Range("I94").Value = 2300
$ws.Range("J94").Value = 64750
$ws.Range("K94").Value = 2300
$ws.Range("L94").Value = 64750
$ws.Range("M94").Value = -1399
$ws.Range("N94").Value = -66552

$ws.Range("H96").Value = 10534655
$ws.Range("J96").Value = 17557192
$ws.Range("L96").Value = 17557192
$ws.Range("N96").Value = -17559938

$ws.Range("H132").Value = 1950
$ws.Range("I132").Value = 1454.6364
$ws.Range("J132").Value = 3312.25
$ws.Range("K132").Value = 4363.9092
$ws.Range("L132").Value = 9936.75
$ws.Range("M132").Value = -1833.9092
$ws.Range("N132").Value = -14996.75

$ws.Range("H136").Value = 948.93335
$ws.Range("I136").Value = 460.30768
$ws.Range("J136").Value = 4125
$ws.Range("K136").Value = 1380.92304
$ws.Range("L136").Value = 12375
$ws.Range("M136").Value = 1169.07696
$ws.Range("N136").Value = -17475
